$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("First sheet")

$values = @{
    "A7" = "Kamran Bains"
    "B7" = "Chloe-Ann Vega"
    "C7" = "Amayah Barajas"
    "D7" = "Safa Blackburn"
    "E7" = "Kezia Gonzalez"
    "F7" = "Boyd Mcbride"
    "G7" = "Leela Romero"
    "H7" = "Mateusz Thornton"
    "I7" = "Amelie Bell"
    "J7" = "Jevon Myers"
    "A8" = "Riley-James Duran"
    "B8" = "Glen Churchill"
    "C8" = "Sachin Deacon"
    "D8" = "Rufus Redfern"
    "E8" = "Jonah Best"
    "F8" = "Zion Ingram"
    "G8" = "Matei Gibbs"
    "H8" = "Kaelan Mcdonnell"
    "I8" = "Spike Peel"
    "J8" = "Zakariyah Gray"
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
